$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the old "_GoBack" bookmark that sat after the "Realizado"
#    run in the "Validar passwords..." bullet.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Append "  Realizado (excepto la ciudad)" to the
#    "Buscar artesanos x cedula, ciudad, por nombre" bullet, with
#    "Realizado" highlighted green + underlined, and
#    " (excepto la ciudad)" underlined only.
# ------------------------------------------------------------------
$target = $d.Content
$target.Find.ClearFormatting()
$target.Find.Execute("Buscar artesanos x cedula, ciudad, por nombre", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target.Collapse(0)

$target.InsertAfter("  ")
$afterSpaces = $target.End

$realizadoRange = $d.Range($afterSpaces, $afterSpaces)
$realizadoRange.InsertAfter("Realizado")
$afterRealizado = $realizadoRange.End

$restRange = $d.Range($afterRealizado, $afterRealizado)
$restRange.InsertAfter(" (excepto la ciudad)")
$afterRest = $restRange.End

$formatRealizado = $d.Range($afterSpaces, $afterRealizado)
$fr = $formatRealizado.Find
$fr.ClearFormatting()
$fr.Replacement.ClearFormatting()
$fr.Replacement.Font.HighlightColorIndex = 4
$fr.Replacement.Font.Underline = 1
$fr.Execute("Realizado", $false, $false, $false, $false, $false, $true, 1, $false, "Realizado", 2) | Out-Null

$formatRest = $d.Range($afterRealizado, $afterRest)
$fr2 = $formatRest.Find
$fr2.ClearFormatting()
$fr2.Replacement.ClearFormatting()
$fr2.Replacement.Font.Underline = 1
$fr2.Execute(" (excepto la ciudad)", $false, $false, $false, $false, $false, $true, 1, $false, " (excepto la ciudad)", 2) | Out-Null

# ------------------------------------------------------------------
# 3. Add the "_GoBack" bookmark back at the end of the very last
#    bullet ("Departamentos y ciudades (para el campo ciudad) y para
#    las búsquedas."), right before its paragraph mark.
# ------------------------------------------------------------------
$lastRange = $d.Content
$lastRange.Find.ClearFormatting()
$lastRange.Find.Execute("Departamentos y ciudades (para el campo ciudad) y para las búsquedas.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$lastRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $lastRange)
